$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are plain text that happen to look numeric
# (dot-grouped thousands, e.g. "69.768.01"). Force text format before
# assigning so Excel does not reinterpret them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.768.01"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.836.12"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.26"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.46"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.834.86"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.87"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.468.84"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.824.99"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.834.65"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.57"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.70"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "508.42"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.745"
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.45"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000144"
$ws.Range("E26").Value = "  +4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.68"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.49"
$ws.Range("E28").Value = "  -6.34%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.53"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.96"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "485.23"
$ws.Range("E39").Value = "  +14.33%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +7.25%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.95"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.57"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.935.69"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.40"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.03"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("E51").Value = "  -2.90%  "
